$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A3").Value = "Mohamed A Talaat"
$ws.Range("B3").Value = "2025-04-17 12:33:02"
